$wb = $excel.ActiveWorkbook

# --- Sheet "OrderClientRef": add a new order/client pairing ---
$wsRef = $wb.Worksheets.Item("OrderClientRef")
$wsRef.Range("A4").Value = 3
$wsRef.Range("B4").Value = 3
[void]$wsRef.Range("E9").Select()

# --- Sheet "OrderList": add the items for the new batch of orders ---
$wsList = $wb.Worksheets.Item("OrderList")

$wsList.Range("A7").Value = 2
$wsList.Range("B7").Value = "hat"
$wsList.Range("C7").Value = 1

$wsList.Range("A8").Value = 3
$wsList.Range("B8").Value = "jacket"
$wsList.Range("C8").Value = 2

$wsList.Range("A9").Value = 3
$wsList.Range("B9").Value = "umbrella"
$wsList.Range("C9").Value = 5

$wsList.Range("A10").Value = 3
$wsList.Range("B10").Value = "rain jacket"
$wsList.Range("C10").Value = 1

$wsList.Range("A11").Value = 3
$wsList.Range("B11").Value = "milk"
$wsList.Range("C11").Value = 1

$wsList.Range("A12").Value = 3
$wsList.Range("B12").Value = "bottle"
$wsList.Range("C12").Value = 10

[void]$wsList.Activate()
[void]$wsList.Range("C14").Select()
